$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.124796628952026
$ws.Range("B1").Value = 2.433515310287476
$ws.Range("C1").Value = 5.168681144714355
$ws.Range("D1").Value = 2.198970794677734
$ws.Range("E1").Value = 1.263920903205872
